$wb = $excel.ActiveWorkbook

# --- 1. Update workbook-level revision/view metadata ---
$wb.RevisionDocumentId = "13_ncr:801_{C83A2452-EFCE-4692-B371-38C670B8E632}"

# --- 2. Add the two new sheets by copying the existing groupStatusCount template ---
$template = $wb.Worksheets.Item("groupStatusCount")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy($null, $lastSheet)
$s4 = $wb.Worksheets.Item($wb.Worksheets.Count)
$s4.Name = "orderTests"

$s4.Copy($null, $s4)
$s5 = $wb.Worksheets.Item($wb.Worksheets.Count)
$s5.Name = "orderTestsStatusHistory"

# --- 3. Fill in orderTests (sheet4) content ---
$s4.Range("A2").Value = "donationId"
$s4.Range("B2").Value = "requestId"
$s4.Range("C2").Value = "Page"
$s4.Range("D2").Value = "Sort"
$s4.Range("E2").Value = "EndPoint"

$s4.Range("A3").Value = "AA0001"
$s4.Range("B3").Value = "AA0001-2020042"
$s4.Range("E3").Value = "/donationInfo/orderTests"

$s4.Range("A6").Value = "donationId"
$s4.Range("B6").Value = "requestId"
$s4.Range("C6").Value = "Page"
$s4.Range("D6").Value = "Sort"
$s4.Range("E6").Value = "EndPoint"

$s4.Range("A7").Value = "AA1"
$s4.Range("B7").Value = "AA0001-2020042"
$s4.Range("E7").Value = "/donationInfo/orderTests"

$s4.Range("A10").Value = "donationId"
$s4.Range("B10").Value = "requestId"
$s4.Range("C10").Value = "Page"
$s4.Range("D10").Value = "Sort"
$s4.Range("E10").Value = "EndPoint"

$s4.Range("A11").Value = "AA0001"
$s4.Range("B11").Value = "AA0001-2020042"
$s4.Range("E11").Value = "/donationInfo/orderTests"

$s4.Rows.Item(2).RowHeight = 30
$s4.Rows.Item(6).RowHeight = 30
$s4.Rows.Item(7).RowHeight = 45
$s4.Rows.Item(10).RowHeight = 30
$s4.Rows.Item(11).RowHeight = 30

$s4.Range("A1").Select()
$s4.Range("A1:E11").Select()

# --- 4. Fill in orderTestsStatusHistory (sheet5) content ---
$s5.Range("A2").Value = "donationId"
$s5.Range("B2").Value = "requestId"
$s5.Range("C2").Value = "Page"
$s5.Range("D2").Value = "Sort"
$s5.Range("E2").Value = "EndPoint"

$s5.Range("A3").Value = "AA0001"
$s5.Range("B3").Value = "AA0001-2020042"
$s5.Range("E3").Value = "/donationInfo/orderTestStatusHistory"

$s5.Range("A6").Value = "donationId"
$s5.Range("B6").Value = "requestId"
$s5.Range("C6").Value = "Page"
$s5.Range("D6").Value = "Sort"
$s5.Range("E6").Value = "EndPoint"

$s5.Range("A7").Value = "AA1"
$s5.Range("B7").Value = "AA0001-2020042"
$s5.Range("E7").Value = "/donationInfo/orderTestStatusHistory"

$s5.Range("A10").Value = "donationId"
$s5.Range("B10").Value = "requestId"
$s5.Range("C10").Value = "Page"
$s5.Range("D10").Value = "Sort"
$s5.Range("E10").Value = "EndPoint"

$s5.Range("A11").Value = "AA0001"
$s5.Range("B11").Value = "AA0001-2020042"
$s5.Range("E11").Value = "/donationInfo/orderTestStatusHistory"

$s5.Rows.Item(2).RowHeight = 30
$s5.Rows.Item(6).RowHeight = 30
$s5.Rows.Item(7).RowHeight = 30
$s5.Rows.Item(10).RowHeight = 30
$s5.Rows.Item(11).RowHeight = 30

$s5.Range("F6").Select()

# --- 5. Tweak existing sheet selections ---
$donationShortDetail = $wb.Worksheets.Item("donationShortDetail")
$donationShortDetail.Range("A1:E11").Select()
$donationShortDetail.Application.ActiveWindow.RangeSelection.Item(1).Activate()

$groupStatusCount = $wb.Worksheets.Item("groupStatusCount")
$groupStatusCount.Range("C11").Select()

# --- 6. Window / view properties ---
$wb.Windows.Item(1).WindowWidth = 12330
$wb.Windows.Item(1).WindowHeight = 9585
$wb.Windows.Item(1).DisplayWorkbookTabs = $true

$s5.Activate()

Write-Output "ok"
